$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's columns were originally laid out as:
#   D = category-code, E = group-code, F = group-name, G = category-name
# The corrected layout (matching the upstream codeforIATI codelist export) is:
#   D = category-code, E = category-name, F = group-code, G = group-name
#
# For every row (including the header row) this is achieved by rotating the
# values held in columns E, F, G one step to the "left":
#   newE = oldG
#   newF = oldE
#   newG = oldF
# Column D (and A, B, C) are left untouched.

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rng = $ws.Range("E1:G$lastRow")

# Every cell in this sheet is stored as text in the source workbook (including
# numeric-looking codes such as "110"). Force the number format to Text before
# writing the rotated values back so Excel does not silently re-type them as
# numbers.
$rng.NumberFormat = "@"

$vals = $rng.Value()

for ($i = 1; $i -le $lastRow; $i++) {
    $oldE = $vals[$i, 1]
    $oldF = $vals[$i, 2]
    $oldG = $vals[$i, 3]

    $vals[$i, 1] = $oldG
    $vals[$i, 2] = $oldE
    $vals[$i, 3] = $oldF
}

$rng.Value = $vals
